$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '69.648.58'
$ws.Cells.Item(2, 5).Value = '  +0.35%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.499.55'
$ws.Cells.Item(3, 5).Value = '  +0.32%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.14%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '603.62'
$ws.Cells.Item(5, 5).Value = '  -0.97%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '195.08'
$ws.Cells.Item(6, 5).Value = '  +4.93%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.29%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.00%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -5.46%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.649'
$ws.Cells.Item(10, 5).Value = '  +0.58%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.00%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000301'
$ws.Cells.Item(12, 5).Value = '  -1.79%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.11%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '4.063.77'
$ws.Cells.Item(14, 5).Value = '  +0.74%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '592.87'
$ws.Cells.Item(15, 5).Value = '  -1.16%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '69.844.76'
$ws.Cells.Item(16, 5).Value = '  +0.53%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'Uniswap'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '12.76'
$ws.Cells.Item(17, 5).Value = '  +1.57%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +1.04%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'WrappedEther'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.507.03'
$ws.Cells.Item(19, 5).Value = '  +0.16%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'TRON'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.122'
$ws.Cells.Item(20, 5).Value = '  +1.74%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.47%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '18.08'
$ws.Cells.Item(22, 5).Value = '  +5.27%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.30'
$ws.Cells.Item(23, 5).Value = '  +4.04%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '102.23'
$ws.Cells.Item(24, 5).Value = '  -3.03%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.55%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +3.39%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.84'
$ws.Cells.Item(27, 5).Value = '  -0.78%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.53'
$ws.Cells.Item(28, 5).Value = '  -1.30%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '33.24'
$ws.Cells.Item(29, 5).Value = '  -0.35%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'dogwifhat'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '4.28'
$ws.Cells.Item(30, 5).Value = '  +3.54%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'NEARProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.03'
$ws.Cells.Item(31, 5).Value = '  +1.44%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '12.41'
$ws.Cells.Item(32, 5).Value = '  +0.17%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.14%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '63.09'
$ws.Cells.Item(34, 5).Value = '  -0.31%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0₃0828'
$ws.Cells.Item(35, 5).Value = '  +6.62%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.708.42'
$ws.Cells.Item(36, 5).Value = '  +3.52%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.09'
$ws.Cells.Item(37, 5).Value = '  -2.33%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.999'
$ws.Cells.Item(38, 5).Value = '  +0.11%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -1.18%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -0.93%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '36.46'
$ws.Cells.Item(41, 5).Value = '  -0.67%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -2.47%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '472.48'
$ws.Cells.Item(43, 5).Value = '  -9.15%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -1.87%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.44%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -4.37%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.29'
$ws.Cells.Item(47, 5).Value = '  -1.51%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.01'
$ws.Cells.Item(48, 5).Value = '  +0.33%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.40'
$ws.Cells.Item(49, 5).Value = '  -4.11%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +1.91%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +9.81%  '
